$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "269.59"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.93"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.390"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06244"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.648"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.701"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.377"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8345"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01377"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1633"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08408"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03406"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03105"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09333"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.890"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001704"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04813"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006339"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001088"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003327"
$ws.Range("E21").Value = "20HotbitTokenHTBWorstin24h"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.734"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.372"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3405"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1264"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04684"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1172"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003438"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01122"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006277"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07551"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002300"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01240"
